$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ShopProductTable")
$v = $ws.Range("A96").Value
Write-Output $v
